$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated exchange rate date (Excel serial date) and last-updated epoch timestamp
$newDate = 44536
$newTime = 1638748801

# Row 2 - AED (rate unchanged)
$ws.Range("D2").Value = $newDate
$ws.Range("E2").Value = $newTime

# Row 3 - AFN
$ws.Range("D3").Value = $newDate
$ws.Range("E3").Value = $newTime
$ws.Range("F3").Value = 96.25

# Row 4 - ALL
$ws.Range("D4").Value = $newDate
$ws.Range("E4").Value = $newTime
$ws.Range("F4").Value = 107.26

# Row 5 - AMD
$ws.Range("D5").Value = $newDate
$ws.Range("E5").Value = $newTime
$ws.Range("F5").Value = 489.44

# Row 6 - ANG (rate unchanged)
$ws.Range("D6").Value = $newDate
$ws.Range("E6").Value = $newTime
